# Regenerate the "K" column (column G) values in the save_data sheet.
# This mirrors a re-run of the underlying calc that produced the K / std / mean
# / s_vals figures, writing the newly computed K values back over the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value
$kValues = @{
    2  = 2
    3  = 3
    4  = 3
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    10 = 2
    11 = 1
    13 = 1
    14 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
